$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.215.06'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.854.95'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7001'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.02'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3088'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07721'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '23.75'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07807'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('D12').Value = '1.865.68'
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '92.01'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.089'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6858'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.32%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.479'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('E17').Value = '  +1.42%  '
$ws.Range('D18').Value = '29.208.76'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '249.27'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '2.112.46'
$ws.Range('E20').Value = '  -3.94%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.83'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.40%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.514'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1528'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '160.11'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.829'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.560'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.06%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.229'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.202'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.193'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.23%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05186'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.33%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7582'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.840'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.711'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01858'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').Value = '1.227.66'
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.725'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.73%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8969'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '109.63'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.551'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -11.30%  '
$ws.Range('D45').Value = '2.009.37'
$ws.Range('E45').Value = '  -3.28%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '65.20'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -9.35%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5183'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.489'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('E49').Value = '  -6.86%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.744'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.990'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.26%  '
